# Add a new course entry ("Ciencia de datos") as a new description row
# above the first table header (row 8), pushing the existing rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 8 (everything from row 8 down shifts to row 9+)
$ws.Rows.Item(8).Insert()

# The new row should look like the other "section description" rows
# (e.g. row 2 / row 4 / row 6), which carry styles s=8 (col A), s=3 (col B,
# merged-look text), s=2 (cols C:J) and s=9 (col K). Copy that formatting
# from row 2 into the freshly inserted row 8.
$ws.Range("A2:K2").Copy()
$ws.Range("A8:K8").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false | Out-Null

# This particular banner is shorter than the others, so it only needs a
# row height of 28 points (the others use the default 42pt wrap height).
$ws.Rows.Item(8).RowHeight = 28

# Fill in the new banner text (wraps to two lines: title + URL).
$ws.Range("B8").Value = "Ciencia de datos`nhttps://cursosacademiaba.buenosaires.gob.ar/"

# Update the sheet's saved selection / view to match what was left selected
# after the edit (no more frozen/scrolled topLeftCell, and the cells that
# were being edited, A13:A14, are selected).
$ws.Range("A13:A14").Select() | Out-Null

Write-Host "Inserted new course row and updated selection"
